$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsBAEP = $wb.Worksheets.Item("BAEPAbCiPC")

# --- About sheet: add California / date columns ---
$wsAbout.Range("B1").Value = "California"
$wsAbout.Range("C1").Value = (Get-Date -Year 2023 -Month 12 -Day 11 -Hour 0 -Minute 0 -Second 0)

# --- BAEPAbCiPC sheet: clear the "electricity/heat/hydrogen" highlight style ---
$wsBAEP.Range("A2:B2").ClearFormats()
$wsBAEP.Range("A15:B15").ClearFormats()
$wsBAEP.Range("A22:B22").ClearFormats()

# --- BAEPAbCiPC sheet: flip values 1 -> 0 for the listed fuel rows ---
$wsBAEP.Range("B3").Value = 0
$wsBAEP.Range("B4").Value = 0
$wsBAEP.Range("B9").Value = 0
$wsBAEP.Range("B10").Value = 0
$wsBAEP.Range("B11").Value = 0
$wsBAEP.Range("B12").Value = 0
$wsBAEP.Range("B13").Value = 0
$wsBAEP.Range("B14").Value = 0
$wsBAEP.Range("B15").Value = 0
$wsBAEP.Range("B17").Value = 0
$wsBAEP.Range("B18").Value = 0
$wsBAEP.Range("B19").Value = 0
$wsBAEP.Range("B20").Value = 0

# --- Selection / view changes ---
$wsBAEP.Range("B18").Select()
$wsAbout.Activate()
$wsBAEP.Activate()
$wb.Windows.Item(1).ScrollRow = 1
